# Team eXi Bug Metrics - iteration 5 update
# Fills in outcome data for the Update Grade bugs (iteration 4, rows 12-14),
# logs two new bugs discovered/fixed in iteration 5 (rows 15-16), and rolls
# the new totals up into the "Bug Metrics" summary sheet (row 12 / H12).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Bug Log" sheet
# ---------------------------------------------------------------------
$log = $wb.Worksheets.Item("Bug Log")

# Row 12 - Wrong label for Assessment Type dropdown: fill in commit id,
# date solved, time taken and status (resolved).
$log.Range("C12").Value = "123d5e72a28deadf1651bc472a7a5fd03c954537"
$log.Range("I12").Value = "19/08/2018"
$log.Range("J12").Value = 1
$log.Range("N12").Value = "Resolved"

# Row 13 - Missing error message when grade update is unsuccessful (front
# end message). Correct the date found and fill in the resolution info.
$log.Range("C13").Value = "123d5e72a28deadf1651bc472a7a5fd03c954537"
$log.Range("G13").Value = "16/08/2018"
$log.Range("I13").Value = "19/08/2018"
$log.Range("J13").Value = 1
$log.Range("N13").Value = "Resolved"
$log.Range("O13").Value = "Added error message in front end"

# Row 14 - Missing validation for the same bug. Correct the date found and
# fill in the resolution info.
$log.Range("C14").Value = "123d5e72a28deadf1651bc472a7a5fd03c954537"
$log.Range("G14").Value = "16/08/2018"
$log.Range("I14").Value = "19/08/2018"
$log.Range("J14").Value = 1
$log.Range("N14").Value = "Resolved"
$log.Range("O14").Value = "Added error message in front end"

# Row 15 - new bug found during iteration 5: Attendance Taking checkbox
# does not disappear once attendance is marked.
$log.Range("B15").Value = 5
$log.Range("D15").Value = "Attendance Taking"
$log.Range("E15").Value = "When attendance is marked, the checkbox doesn’t disappear and users are confused"
$log.Range("F15").Value = "Poor Interface Design"
$log.Range("G15").Value = "31/08/2018"
$log.Range("H15").Value = "Hui Xin"
$log.Range("I15").Value = "31/08/2018"
$log.Range("K15").Value = "Low"
$log.Range("L15").Value = 1
$log.Range("M15").Value = 38

# Row 16 - new bug found during iteration 5: Manage Acc page crashes
# because a renamed DAO method wasn't updated everywhere.
$log.Range("B16").Value = 5
$log.Range("D16").Value = "Manage Acc"
$log.Range("E16").Value = "Crashes when accessing page"
$log.Range("F16").Value = "Did not change the method name when we renamed it in the DAO"
$log.Range("G16").Value = "31/08/2019"
$log.Range("H16").Value = "Hui Xin"
$log.Range("I16").Value = "31/08/2018"
$log.Range("J16").Value = 0.1
$log.Range("K16").Value = "Critical"
$log.Range("L16").Value = 10
$log.Range("M16").Value = 48
$log.Range("N16").Value = "Resolved"
$log.Range("O16").Value = "Changed the method name accordingly"
$log.Range("C16").Value = "669ddc1f9a6e3d8c79a4a2a85ba046489fe0d000"

$log.Range("C16").Select()

# ---------------------------------------------------------------------
# "Bug Metrics" sheet
# ---------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Bug Metrics")

# Iteration 5 summary row.
$metrics.Range("C12").Value = 2
$metrics.Range("D12").Value = 1
$metrics.Range("E12").Value = 0
$metrics.Range("F12").Value = 1
$metrics.Range("G12").Value = 48
$metrics.Range("H12").Value = "The bugs were discovered during the testing of the functionalities and hence we fixed it immediately upon discovery"

$metrics.Range("H12:K12").Select()
